# Add columns I (I0) and J (IF) to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) — copy the formatting of the existing header cell (H1,
# bold/bordered/centered style) onto the two new header cells, then set
# their text.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# Data values for rows 2..63
$iVals = @(8,8,9,9,7,7,6,6,8,9,7,7,8,7,10,5,9,9,7,7,7,7,8,9,8,10,8,7,7,7,6,6,6,7,7,7,7,8,10,7,7,7,9,6,6,8,7,7,7,7,8,6,7,7,8,7,7,8,6,6,9,8)
$jVals = @(8,8,9,9,7,8,6,7,8,9,7,8,8,8,10,6,9,9,7,7,7,7,8,9,8,10,8,8,7,7,6,6,6,8,8,7,8,8,10,7,7,8,9,6,6,8,7,7,7,7,8,6,7,7,8,7,7,8,6,6,9,8)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
